# Update the East-Asian / complex-script fonts used across the document's
# paragraph styles (vignettes docx/html regeneration - "update docx and
# html under vignettes"):
#   - docDefaults / Normal / Heading: eastAsia font DejaVu Sans -> Tahoma
#   - List / Caption / Index: gain an explicit complex-script (cs) font of
#     DejaVu Sans, matching what they already inherited.

$d = $word.ActiveDocument
$styles = $d.Styles

# Normal: ascii/hAnsi = Liberation Serif, eastAsia DejaVu Sans -> Tahoma, cs stays DejaVu Sans
$normal = $styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

# Heading: ascii/hAnsi = Liberation Sans, eastAsia DejaVu Sans -> Tahoma, cs stays DejaVu Sans
$heading = $styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

# List: previously empty rPr, now explicit complex-script font DejaVu Sans
$list = $styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

# Caption: rPr gains a complex-script font DejaVu Sans ahead of the existing italics
$caption = $styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

# Index: previously empty rPr, now explicit complex-script font DejaVu Sans
$index = $styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
